# Update "想去人数" (column F) values across sheets per the scraped snapshot
# at commit 456a3b4 (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 19290
$ws.Cells.Item(5, 6).Value = 767
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 1084
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 7308
$ws.Cells.Item(10, 6).Value = 470
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(14, 6).Value = 141
$ws.Cells.Item(15, 6).Value = 90
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(31, 6).Value = 547
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 3
$ws.Cells.Item(41, 6).Value = 246
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 3962
$ws.Cells.Item(45, 6).Value = 0

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 2

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 26
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 228
$ws.Cells.Item(18, 6).Value = 1322
$ws.Cells.Item(20, 6).Value = 65
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 152
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 6).Value = 2
$ws.Cells.Item(33, 6).Value = 39
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 23
$ws.Cells.Item(37, 6).Value = 82
$ws.Cells.Item(38, 6).Value = 12402
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 51
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(44, 6).Value = 325
$ws.Cells.Item(45, 6).Value = 0

Write-Output "Updated F-column values on 展览, 演出, 全部类型"
